$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the first 8 data rows (16-23) so the remaining 9 rows (old 24-32)
# shift up to become rows 16-24. This preserves the special "closing" border
# style that lives on the old last row (32) by landing it on the new last
# row (24), while the plain style on old rows 24-31 lands on new rows 16-23.
$ws.Rows("16:23").Delete()

# --- Update summary header cells ---
$ws.Range("E11").Value = 376844      # VALOR MORA total
$ws.Range("C13").Value = 5           # Cant. Trabajadores
$ws.Range("F13").Value = 7           # Cant. Periodos

# --- Overwrite the 9 remaining data rows (16-24) with the new dataset ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1065377933"
$ws.Range("D16").Value = "DARIO JOSE OROZCO CARVAJAL"
$ws.Range("E16").Value = "2103"
$ws.Range("F16").Value = 48000
$ws.Range("G16").Value = 1519000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1065377933"
$ws.Range("D17").Value = "DARIO JOSE OROZCO CARVAJAL"
$ws.Range("E17").Value = "2104"
$ws.Range("F17").Value = 48000
$ws.Range("G17").Value = 1519000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "9297885"
$ws.Range("D18").Value = "CAMILO ENRIQUE ZABALETA CARDONA"
$ws.Range("E18").Value = "2206"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1519000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1007739301"
$ws.Range("D19").Value = "JONATAN MICHAEL FORERO AVILA"
$ws.Range("E19").Value = "2307"
$ws.Range("F19").Value = 46400
$ws.Range("G19").Value = 1160000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1007739301"
$ws.Range("D20").Value = "JONATAN MICHAEL FORERO AVILA"
$ws.Range("E20").Value = "2308"
$ws.Range("F20").Value = 46400
$ws.Range("G20").Value = 1160000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45516414"
$ws.Range("D21").Value = "PATRICIA DEL ROSARIO ARZUZA DIAZ"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 17082
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1094161333"
$ws.Range("D22").Value = "CLEIVER ASET BAYONA RUIZ"
$ws.Range("E22").Value = "2507"
$ws.Range("F22").Value = 17082
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "45516414"
$ws.Range("D23").Value = "PATRICIA DEL ROSARIO ARZUZA DIAZ"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1094161333"
$ws.Range("D24").Value = "CLEIVER ASET BAYONA RUIZ"
$ws.Range("E24").Value = "2508"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

# --- Column D was auto-fit to the new (narrower) longest name ---
$ws.Columns("D").AutoFit() | Out-Null
